$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new shared strings in the exact order they must be appended
# to xl/sharedStrings.xml (index 387..410), which is NOT the same as the
# natural row/column traversal order of the new table rows (158..165).
$ws.Range("B158").Value = ' : Line {0} Col {1}'
$ws.Range("C158").Value = ' : 行 {0} 列 {1}'
$ws.Range("A159").Value = 'CompileErrror_InvalidToken'
$ws.Range("A160").Value = 'CompileErrror_InvalidEoF'
$ws.Range("B159").Value = 'Invalid token {0} is found'
$ws.Range("B160").Value = 'Unexpected end is found'
$ws.Range("C160").Value = '予期せぬ終了が見つかりました'
$ws.Range("C159").Value = '無効なトークン {0} が見つかりました'
$ws.Range("A161").Value = 'CompileError_InvalidOperation'
$ws.Range("B161").Value = 'Invalid operation is found'
$ws.Range("C161").Value = '無効な操作が見つかりました'
$ws.Range("A162").Value = 'CompileError_ArgSize'
$ws.Range("C162").Value = '引数の個数が間違っています。要求 {0} 実際 {1}'
$ws.Range("B162").Value = 'The number of arguments is wrong. Expected {0} Actual {1}'
$ws.Range("B163").Value = 'Invalid substitution is found'
$ws.Range("C163").Value = '無効な代入が見つかりました'
$ws.Range("A163").Value = 'CompileError_InvalidSubstitution'
$ws.Range("B164").Value = 'Unkwon function {0} is found'
$ws.Range("C164").Value = '知らない関数 {0} が見つかりました'
$ws.Range("A164").Value = 'CompileError_UnknownFunction'
$ws.Range("A165").Value = 'CompileError_UnknownValue'
$ws.Range("B165").Value = 'Unkwon value {0} is found'
$ws.Range("C165").Value = '知らない値 {0} が見つかりました'
$ws.Range("A158").Value = 'CompilePosition'

# Column A of the first new row (158) carries the same style as the rest
# of the table (vertical-centered, wrapped text) -- matches style index 1
# used throughout the rest of column A/B/C in this sheet.
$ws.Range("A158").WrapText = $true
$ws.Range("A158").VerticalAlignment = -4108

# Update the visible selection to match the post-edit state recorded in
# the workbook (user had just finished typing the new key column).
$ws.Range("A159:A165").Select() | Out-Null
